$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "15979240"
$ws.Range("H2").Value = "080-01-6913172"

$ws.Range("F10").Select()
